# Update Daily Report: 2026-02-24
# Appends the 2026-02-23 (serial 46076) daily rows to Daily_Data, and
# refreshes the dependent Today_Summary / Monthly_Stats roll-ups.

$wb = $excel.ActiveWorkbook

$wsDaily   = $wb.Worksheets.Item("Daily_Data")
$wsSummary = $wb.Worksheets.Item("Today_Summary")
$wsMonthly = $wb.Worksheets.Item("Monthly_Stats")

# ---------------------------------------------------------------------------
# 1. Daily_Data: append rows 242-265 (date serial 46076 / 2026-02-23)
#    Columns: Date, Region_Type, PREV_TOTAL, RECEIVED, WITHDRAWN, NET_CHANGE,
#             ADJUSTMENT, TOTAL_TODAY
# ---------------------------------------------------------------------------
$newRows = @(
    @(242, 46076, 'ASAHI DEPOSITORY LLC Registered', 23301775.992, 0, 0, 0, -10160, 23291615.992),
    @(243, 46076, 'ASAHI DEPOSITORY LLC Eligible', 2748893.808, 0, 0, 0, 10160, 2759053.808),
    @(244, 46076, "BRINK'S, INC. Registered", 15782712.636, 0, 0, 0, -20279.63, 15762433.006),
    @(245, 46076, "BRINK'S, INC. Eligible", 39336942.517, 0, 0, 0, 20279.63, 39357222.147),
    @(246, 46076, 'CNT DEPOSITORY, INC. Registered', 12174851.569, 0, 0, 0, 0, 12174851.569),
    @(247, 46076, 'CNT DEPOSITORY, INC. Eligible', 13856687.823, 0, 0, 0, 0, 13856687.823),
    @(248, 46076, 'DELAWARE DEPOSITORY Registered', 1532776.423, 0, 0, 0, 0, 1532776.423),
    @(249, 46076, 'DELAWARE DEPOSITORY Eligible', 15771876.945, 0, 2968.9, -2968.9, 0, 15768908.045),
    @(250, 46076, 'HSBC BANK, USA Registered', 3412157.57, 0, 0, 0, -19930.44, 3392227.13),
    @(251, 46076, 'HSBC BANK, USA Eligible', 19011209.353, 0, 0, 0, 19930.44, 19031139.793),
    @(252, 46076, 'INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Registered', 273789.87, 0, 0, 0, 0, 273789.87),
    @(253, 46076, 'INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Eligible', 3295246.644, 0, 0, 0, 0, 3295246.644),
    @(254, 46076, 'JP MORGAN CHASE BANK NA Registered', 12000343.77, 0, 0, 0, 0, 12000343.77),
    @(255, 46076, 'JP MORGAN CHASE BANK NA Eligible', 142890665.333, 0, 0, 0, 0, 142890665.333),
    @(256, 46076, 'LOOMIS INTERNATIONAL (US) LLC Registered', 6311885.937, 0, 0, 0, -5046.5, 6306839.437),
    @(257, 46076, 'LOOMIS INTERNATIONAL (US) LLC Eligible', 24033585.186, 0, 0, 0, 5046.5, 24038631.686),
    @(258, 46076, 'MALCA-AMIT ARMORED, INC. Registered', 0, 0, 0, 0, 0, 0),
    @(259, 46076, 'MALCA-AMIT ARMORED, INC. Eligible', 0, 0, 0, 0, 0, 0),
    @(260, 46076, 'MALCA-AMIT USA, LLC Registered', 949634.064, 0, 0, 0, 0, 949634.064),
    @(261, 46076, 'MALCA-AMIT USA, LLC Eligible', 1073898.377, 0, 0, 0, 0, 1073898.377),
    @(262, 46076, 'MANFRA, TORDELLA & BROOKES, LLC Registered', 6219630.033, 0, 0, 0, 0, 6219630.033),
    @(263, 46076, 'MANFRA, TORDELLA & BROOKES, LLC Eligible', 12256015.907, 0, 0, 0, 0, 12256015.907),
    @(264, 46076, 'STONEX PRECIOUS METALS LLC Registered', 6231501.4, 0, 0, 0, -963345.72, 5268155.68),
    @(265, 46076, 'STONEX PRECIOUS METALS LLC Eligible', 1537051.72, 0, 0, 0, 963345.72, 2500397.44)
)

foreach ($row in $newRows) {
    $r = $row[0]
    $dateCell = $wsDaily.Cells.Item($r, 1)
    $dateCell.Value = $row[1]
    $dateCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $wsDaily.Cells.Item($r, 2).Value = $row[2]
    $wsDaily.Cells.Item($r, 3).Value = $row[3]
    $wsDaily.Cells.Item($r, 4).Value = $row[4]
    $wsDaily.Cells.Item($r, 5).Value = $row[5]
    $wsDaily.Cells.Item($r, 6).Value = $row[6]
    $wsDaily.Cells.Item($r, 7).Value = $row[7]
    $wsDaily.Cells.Item($r, 8).Value = $row[8]
}

# ---------------------------------------------------------------------------
# 2. Today_Summary: refresh Eligible / Registered (and derived Total_Stock)
#    for the depositories touched by the new day's adjustments/withdrawals.
# ---------------------------------------------------------------------------
$summaryUpdates = @(
    @("ASAHI DEPOSITORY LLC", 2759053.808, 23291615.992, $false),
    @("BRINK'S, INC.", 39357222.147, 15762433.006, $false),
    @("DELAWARE DEPOSITORY", 15768908.045, 1532776.423, $true),
    @("HSBC BANK, USA", 19031139.793, 3392227.13, $false),
    @("LOOMIS INTERNATIONAL (US) LLC", 24038631.686, 6306839.437, $false),
    @("STONEX PRECIOUS METALS LLC", 2500397.44, 5268155.68, $true)
)

$lastSummaryRow = $wsSummary.UsedRange.Rows.Count
for ($i = 2; $i -le $lastSummaryRow; $i++) {
    $name = $wsSummary.Cells.Item($i, 1).Value2
    foreach ($u in $summaryUpdates) {
        if ($name -eq $u[0]) {
            $wsSummary.Cells.Item($i, 2).Value = $u[1]
            $wsSummary.Cells.Item($i, 3).Value = $u[2]
            if ($u[3]) {
                $wsSummary.Cells.Item($i, 4).Value = $u[1] + $u[2]
            }
        }
    }
}

# ---------------------------------------------------------------------------
# 3. Monthly_Stats: refresh the month overview (row 2) and the per
#    depository/type detail rows (7-30) for February 2026.
# ---------------------------------------------------------------------------
$monthlyDetailUpdates = @(
    @("ASAHI DEPOSITORY LLC Eligible", 2759053.808, $null),
    @("ASAHI DEPOSITORY LLC Registered", 23291615.992, $null),
    @("BRINK'S, INC. Eligible", 39357222.147, $null),
    @("BRINK'S, INC. Registered", 15762433.006, $null),
    @("DELAWARE DEPOSITORY Eligible", 15768908.045, 694388.0240000001),
    @("HSBC BANK, USA Eligible", 19031139.793, $null),
    @("HSBC BANK, USA Registered", 3392227.13, $null),
    @("LOOMIS INTERNATIONAL (US) LLC Eligible", 24038631.686, $null),
    @("LOOMIS INTERNATIONAL (US) LLC Registered", 6306839.437, $null),
    @("STONEX PRECIOUS METALS LLC Eligible", 2500397.44, $null),
    @("STONEX PRECIOUS METALS LLC Registered", 5268155.68, $null)
)

$lastMonthlyRow = $wsMonthly.UsedRange.Rows.Count
for ($i = 7; $i -le $lastMonthlyRow; $i++) {
    $name = $wsMonthly.Cells.Item($i, 2).Value2
    foreach ($u in $monthlyDetailUpdates) {
        if ($name -eq $u[0]) {
            $wsMonthly.Cells.Item($i, 5).Value = $u[1]
            if ($null -ne $u[2]) {
                $wsMonthly.Cells.Item($i, 4).Value = $u[2]
            }
        }
    }
}

# Refresh the Eligible / Registered / Grand_Total month overview to match
# the resummed per-row detail (rows 7-30) for February 2026.
$wsMonthly.Cells.Item(2, 2).Value = 276827867.003
$wsMonthly.Cells.Item(2, 3).Value = 87172296.97399999
$wsMonthly.Cells.Item(2, 4).Value = 364000163.977
